$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The blank "Thanksgiving" spacer row (row 29) is removed entirely; every
# row below it shifts up by one, and formulas/number formats/fills that
# lived on those rows move up with them automatically.
$ws.Rows(29).Delete() | Out-Null

# Leave the selection where the editor ended up after the delete.
$ws.Range("D35").Select() | Out-Null
